$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 picks up an explicit (formatted) style, matching the rest of
# the data rows above it -- this is what Excel does when a previously
# "raw" row gets touched/re-saved and inherits the sheet's normal style. ---
$ws.Range("A11:D11").Style = "Normal"

# --- Append the new registration entry as row 12. ---
# Format the row as Text first so the numeric-looking phone number
# ("1234567777") is stored as a literal string, like the other phone
# numbers in this sheet, instead of being auto-converted to a number.
$ws.Range("A12:D12").NumberFormat = "@"

$ws.Range("A12").Value = "uhadslkjhkfjh"
$ws.Range("B12").Value = "1234567777"
$ws.Range("C12").Value = "ajhsd@jash.com"
$ws.Range("D12").Value = "asjhd"
